$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item(1)
# Row 6 (hunk 0)
$ws.Range("H6").Value = 408.66666
$ws.Range("I6").Value = 240.4
$ws.Range("J6").Value = 1250
$ws.Range("K6").Value = 721.2
$ws.Range("L6").Value = 3750
$ws.Range("M6").Value = -609.2
$ws.Range("N6").Value = -3974

# Row 49 (hunk 1)
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents() | Out-Null

# Row 70 (hunk 2)
$ws.Range("H70").Value = 2450
$ws.Range("I70").Value = 1800
$ws.Range("J70").Value = 2666.6667
$ws.Range("K70").Value = 5400
$ws.Range("L70").Value = 8000.000100000001
$ws.Range("M70").Value = -5130
$ws.Range("N70").Value = -8540.000100000001

# Row 73 (hunk 3)
$ws.Range("H73").Value = 2450
$ws.Range("I73").Value = 1800
$ws.Range("J73").Value = 2666.6667
$ws.Range("K73").Value = 5400
$ws.Range("L73").Value = 8000.000100000001
$ws.Range("M73").Value = -4464
$ws.Range("N73").Value = -9872.000100000001

# Row 132 (hunk 4)
$ws.Range("H132").Value = 37222.16
$ws.Range("I132").Value = 37567.215
$ws.Range("K132").Value = 112701.645
$ws.Range("M132").Value = -110171.645

# Sheet: ARM
$ws = $wb.Worksheets.Item(2)
# Row 49 (hunk 5)
$ws.Range("H49").Value = 5300
$ws.Range("J49").Value = 5300
$ws.Range("L49").Value = 5300
$ws.Range("N49").Value = -5820

# Row 114 (hunk 6)
$ws.Range("H114").Value = 32000
$ws.Range("J114").Value = 32000
$ws.Range("L114").Value = 32000
$ws.Range("N114").Value = -40678

# Row 122 (hunk 7)
$ws.Range("H122").Value = 3463.75
$ws.Range("I122").Value = 2674.56
$ws.Range("J122").Value = 6282.2856
$ws.Range("K122").Value = 8023.68
$ws.Range("L122").Value = 18846.8568
$ws.Range("M122").Value = -5573.68
$ws.Range("N122").Value = -23746.8568

# Sheet: BSM
$ws = $wb.Worksheets.Item(3)
# Row 99 (hunk 8)
$ws.Range("H99").Value = 1627.7407
$ws.Range("I99").Value = 1606.7727
$ws.Range("K99").Value = 1606.7727
$ws.Range("M99").Value = -108.7727

# Row 100 (hunk 9)
$ws.Range("H100").Value = 21321.5
$ws.Range("J100").Value = 21321.5
$ws.Range("L100").Value = 21321.5
$ws.Range("N100").Value = -23485.5

# Sheet: CRP
$ws = $wb.Worksheets.Item(4)
# Row 31 (hunk 10)
$ws.Range("H31").Value = 5992.636
$ws.Range("I31").Value = 2177.6667
$ws.Range("J31").Value = 9171.777
$ws.Range("K31").Value = 2177.6667
$ws.Range("L31").Value = 9171.777
$ws.Range("M31").Value = -1882.6667
$ws.Range("N31").Value = -9761.777

# Row 34 (hunk 11)
$ws.Range("H34").Value = 5992.636
$ws.Range("I34").Value = 2177.6667
$ws.Range("J34").Value = 9171.777
$ws.Range("K34").Value = 2177.6667
$ws.Range("L34").Value = 9171.777
$ws.Range("M34").Value = -1975.6667
$ws.Range("N34").Value = -9575.777

# Row 39 (hunk 12)
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("M39").ClearContents() | Out-Null
$ws.Range("N39").ClearContents() | Out-Null

# Row 44 (hunk 13)
$ws.Range("H44").Value = 6400
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 6400
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 6400
$ws.Range("M44").ClearContents() | Out-Null
$ws.Range("N44").Value = -7284

# Row 49 (hunk 14)
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").ClearContents() | Out-Null
$ws.Range("N49").ClearContents() | Out-Null

# Row 122 (hunk 15)
$ws.Range("H122").Value = 1301.2307
$ws.Range("I122").Value = 960.2
$ws.Range("J122").Value = 2438
$ws.Range("K122").Value = 2880.6
$ws.Range("L122").Value = 7314
$ws.Range("M122").Value = -430.6000000000004
$ws.Range("N122").Value = -12214

# Row 134 (hunk 16)
$ws.Range("H134").Value = 2931.0386
$ws.Range("I134").Value = 1502.7
$ws.Range("K134").Value = 4508.1
$ws.Range("M134").Value = -1973.1

# Sheet: CUL
$ws = $wb.Worksheets.Item(5)
# Row 47 (hunk 17)
$ws.Range("H47").Value = 534.3333
$ws.Range("I47").Value = 301.5
$ws.Range("K47").Value = 904.5
$ws.Range("M47").Value = -473.5

# Row 49 (hunk 18)
$ws.Range("H49").Value = 743
$ws.Range("I49").Value = 743
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 2229
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -2073
$ws.Range("N49").ClearContents() | Out-Null

# Row 113 (hunk 19)
$ws.Range("H113").Value = 21739820
$ws.Range("I113").Value = 700
$ws.Range("J113").Value = 31250684
$ws.Range("K113").Value = 2100
$ws.Range("L113").Value = 93752052
$ws.Range("M113").Value = 70
$ws.Range("N113").Value = -93756392

# Row 140 (hunk 20)
$ws.Range("H140").Value = 4091.5386
$ws.Range("I140").Value = 4865.4165
$ws.Range("J140").Value = 2853.3333
$ws.Range("K140").Value = 14596.2495
$ws.Range("L140").Value = 8559.999899999999
$ws.Range("M140").Value = -9416.249500000002
$ws.Range("N140").Value = -18919.9999

# Sheet: GSM
$ws = $wb.Worksheets.Item(6)
# Row 20 (hunk 21)
$ws.Range("H20").Value = 7000
$ws.Range("J20").Value = 7000
$ws.Range("L20").Value = 7000
$ws.Range("N20").Value = -7490

# Row 47 (hunk 22)
$ws.Range("H47").Value = 8000
$ws.Range("J47").Value = 8000
$ws.Range("L47").Value = 8000
$ws.Range("N47").Value = -9136

# Row 95 (hunk 23)
$ws.Range("H95").Value = 14958.333
$ws.Range("J95").Value = 14958.333
$ws.Range("L95").Value = 14958.333
$ws.Range("N95").Value = -20450.333

# Row 103 (hunk 24)
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents() | Out-Null

# Row 122 (hunk 25)
$ws.Range("H122").Value = 1814.6428
$ws.Range("I122").Value = 1521
$ws.Range("J122").Value = 2891.3333
$ws.Range("K122").Value = 4563
$ws.Range("L122").Value = 8673.999899999999
$ws.Range("M122").Value = -2113
$ws.Range("N122").Value = -13573.9999

# Sheet: LTW
$ws = $wb.Worksheets.Item(7)
# Row 7 (hunk 26)
$ws.Range("H7").Value = 3430.25
$ws.Range("I7").Value = 3466.6667
$ws.Range("J7").Value = 3423.8235
$ws.Range("K7").Value = 3466.6667
$ws.Range("L7").Value = 3423.8235
$ws.Range("M7").Value = -3354.6667
$ws.Range("N7").Value = -3647.8235

# Row 42 (hunk 27)
$ws.Range("H42").Value = 50000
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents() | Out-Null

# Row 49 (hunk 28)
$ws.Range("H49").Value = 50000
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents() | Out-Null

# Row 51 (hunk 29)
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents() | Out-Null

# Row 126 (hunk 30)
$ws.Range("H126").Value = 3430.25
$ws.Range("I126").Value = 3466.6667
$ws.Range("J126").Value = 3423.8235
$ws.Range("K126").Value = 10400.0001
$ws.Range("L126").Value = 10271.4705
$ws.Range("M126").Value = -7930.000100000001
$ws.Range("N126").Value = -15211.4705

# Sheet: WVR
$ws = $wb.Worksheets.Item(8)
# Row 47 (hunk 31)
$ws.Range("H47").Value = 8237.25
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 8237.25
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 8237.25
$ws.Range("M47").ClearContents() | Out-Null
$ws.Range("N47").Value = -9381.25

# Row 49 (hunk 32)
$ws.Range("H49").Value = 4888.5
$ws.Range("I49").Value = 3777
$ws.Range("K49").Value = 3777
$ws.Range("M49").Value = -3547

# Row 62 (hunk 33)
$ws.Range("H62").Value = 4181600
$ws.Range("I62").Value = 8350525
$ws.Range("J62").Value = 12675.083
$ws.Range("K62").Value = 8350525
$ws.Range("L62").Value = 12675.083
$ws.Range("M62").Value = -8349901
$ws.Range("N62").Value = -13923.083

# Row 65 (hunk 34)
$ws.Range("H65").Value = 4181600
$ws.Range("I65").Value = 8350525
$ws.Range("J65").Value = 12675.083
$ws.Range("K65").Value = 41752625
$ws.Range("L65").Value = 63375.415
$ws.Range("M65").Value = -41749505
$ws.Range("N65").Value = -69615.41500000001

# Row 75 (hunk 35)
$ws.Range("H75").Value = 265059
$ws.Range("J75").Value = 30000
$ws.Range("L75").Value = 30000
$ws.Range("N75").Value = -31872

# Row 78 (hunk 36)
$ws.Range("H78").Value = 265059
$ws.Range("J78").Value = 30000
$ws.Range("L78").Value = 90000
$ws.Range("N78").Value = -99360

# Row 125 (hunk 37)
$ws.Range("H125").Value = 26650
$ws.Range("J125").Value = 26650
$ws.Range("L125").Value = 26650
$ws.Range("N125").Value = -36490
